$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- Row 30 ----
$ws.Range("A30").Value = 112414988
$ws.Range("B30").Value = 90823
$ws.Range("C30").Value = "Ovaliderad"
$ws.Range("D30").Value = "NT"
$ws.Range("E30").Value = 5966
$ws.Range("F30").Value = "Motaggsvamp"
$ws.Range("G30").Value = "Sarcodon squamosus"
$ws.Range("H30").Value = "(Schaeff.) Quél."
$ws.Range("P30").Value = "Spångmyran, Dlr"
$ws.Range("Q30").Value = 381294
$ws.Range("R30").Value = 6862860
$ws.Range("S30").Value = 2
$ws.Range("T30").Value = "Dalarna"
$ws.Range("U30").Value = "Älvdalen"
$ws.Range("V30").Value = "Dalarna"
$ws.Range("W30").Value = "Idre"
$ws.Range("Y30").NumberFormat = "@"
$ws.Range("Y30").Value = "2023-09-30"
$ws.Range("Z30").Value = "11:11"
$ws.Range("AA30").NumberFormat = "@"
$ws.Range("AA30").Value = "2023-09-30"
$ws.Range("AB30").Value = "11:11"
$ws.Range("AD30").Value = $false
$ws.Range("AE30").Value = $false
$ws.Range("AG30").Value = $false
$ws.Range("AW30").Value = "Ingunn Woldmo"
$ws.Range("AX30").Value = "Ingunn Woldmo"

# ---- Row 31 ----
$ws.Range("A31").Value = 112415012
$ws.Range("B31").Value = 90812
$ws.Range("C31").Value = "Ovaliderad"
$ws.Range("D31").Value = "LC"
$ws.Range("E31").Value = 4366
$ws.Range("F31").Value = "Skarp dropptaggsvamp"
$ws.Range("G31").Value = "Hydnellum peckii"
$ws.Range("H31").Value = "Banker"
$ws.Range("P31").Value = "Spångmyran, Dlr"
$ws.Range("Q31").Value = 381335
$ws.Range("R31").Value = 6862894
$ws.Range("S31").Value = 2
$ws.Range("T31").Value = "Dalarna"
$ws.Range("U31").Value = "Älvdalen"
$ws.Range("V31").Value = "Dalarna"
$ws.Range("W31").Value = "Idre"
$ws.Range("Y31").NumberFormat = "@"
$ws.Range("Y31").Value = "2023-09-30"
$ws.Range("Z31").Value = "11:15"
$ws.Range("AA31").NumberFormat = "@"
$ws.Range("AA31").Value = "2023-09-30"
$ws.Range("AB31").Value = "11:15"
$ws.Range("AD31").Value = $false
$ws.Range("AE31").Value = $false
$ws.Range("AG31").Value = $false
$ws.Range("AW31").Value = "Ingunn Woldmo"
$ws.Range("AX31").Value = "Ingunn Woldmo"
